# Added EM data for 13/2/2024
#
# The source workbook keeps one data row per observation date for each of the
# two simulations (ForestHill2023IrrigationFull / ...Partial). A new
# observation date (13/2/2024 -> serial 45335) is being added for BOTH
# simulations:
#   - "ForestHill2023IrrigationFull" already has a placeholder further down
#     the sheet (it is date-sorted per simulation block); the new row is
#     inserted right above the previous first "Full" row (row 19), pushing
#     every following row down by one.
#   - "ForestHill2023IrrigationPartial" gets its new observation appended
#     as the last populated row of its block (new row 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the existing row 19 -----------------------
# (shifts the old rows 19-35 down to 20-36; formatting of the surrounding
# rows is carried onto the new row automatically, same as Excel's own
# "Insert Row" command.)
$ws.Rows.Item(19).Insert()

# --- 2. Populate the newly inserted row 19 (EMp100/EMv100 for the "Full"
#        simulation on 13/2/2024) ----------------------------------------
$ws.Range("A19").Value = "ForestHill2023IrrigationFull"
$ws.Range("B19").Value = 45335
$ws.Range("K19").Value = 191.41666666666666
$ws.Range("K19").Style = "Normal"
$ws.Range("L19").Value = 164.79999999999998

# --- 3. Append the new trailing row (37) for the "Partial" simulation ----
$ws.Range("A37").Value = "ForestHill2023IrrigationPartial"

$ws.Range("B37").Value = 45335
$ws.Range("B37").NumberFormat = "mm-dd-yy"

$ws.Range("K37").Value = 146.33333333333334
$ws.Range("K37").Style = "Normal"
$ws.Range("K37").NumberFormat = "0.00"

$ws.Range("L37").Value = 103.91666666666669
$ws.Range("L37").Style = "Normal"

# --- 4. Extend the cached _FilterDatabase range by the one inserted row --
foreach ($n in $wb.Names) {
    if ($n.Name() -like "*_FilterDatabase*") {
        $n.RefersTo = "=CottonObserved!`$A`$1:`$EQ`$2579"
    }
}

# --- 5. Leave the cursor where the author left it (bottom-right pane on
#        the newly added "Full" observation) ------------------------------
[void]$ws.Range("B37").Select()
[void]$ws.Range("K19").Select()
